$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 18 (the existing "Granada" record) down into a new row 19,
# shifting nothing else since row 19 is currently blank/non-existent.
$ws.Rows.Item(18).Copy()
$ws.Rows.Item(19).Insert()

# Row 18 now becomes the new, updated record (newer date, new volume/prices/origin).
$ws.Range("D18").Value = 44714
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("Q18").Value = "$/caja 18 kilos granel"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 1111
$ws.Range("T18").Value = 18
